# Progress Tracker update: mark several tasks on the "Rubric-I" sheet as
# complete (put an "X" in the "Earned" column) now that ticket/payment
# features have been wired up to the front end.
#
# Rows affected on "Rubric-I" (col E = Earned/"C" marker column):
#   6  - Browse Flights
#   7  - Select Flight
#   8  - Browse Seat Map and Select Seat
#   9  - Make Payment
#   12 - Browse Passenger List by Airline Staff Only
# Each row's F column holds =IF(E=="X", SUM(C:D), 0) so it recalculates
# automatically once E is set.

$wb = $excel.ActiveWorkbook

$rubricI = $wb.Worksheets.Item("Rubric-I")
$tl = $wb.Worksheets.Item("TL")

# Mark the newly-finished tasks as earned.
$rubricI.Range("E6").Value = "X"
$rubricI.Range("E7").Value = "X"
$rubricI.Range("E8").Value = "X"
$rubricI.Range("E9").Value = "X"
$rubricI.Range("E12").Value = "X"

# Restore the "TL" sheet's own selection before moving away from it, so it
# doesn't keep the stale B7 selection.
$tl.Range("B10").Select() | Out-Null

# Bring "Rubric-I" to the front and leave the selection where the author
# left off.
$rubricI.Activate() | Out-Null
$rubricI.Range("E13").Select() | Out-Null
